$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "desc"
$ws.Range("D2").Value = "aaaaa"
$ws.Range("D3").Value = "bbbbb"

$ws.Range("E7").Select()
